$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Insert a brand-new "2022-Q4" sheet right before "2022-Q3"
# ---------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Add($q3)
$ws.Name = "2022-Q4"

# Match the page-margin conventions used by the sibling quarter sheets
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row (bold, centered, thin-bordered - matches the other quarter sheets)
$headers = "基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名"
$col = 2
foreach ($h in $headers) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $col = $col + 1
}

# Data rows
$a = $ws.Cells.Item(2, 1)
$a.Value = 0
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "002446"
$ws.Range("C2").Value = "广发利鑫灵活配置混合A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "73.90"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2.14"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.4821"
$ws.Range("H2").Value = 10

$a = $ws.Cells.Item(3, 1)
$a.Value = 1
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "400032"
$ws.Range("C3").Value = "东方主题精选混合"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "12.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "92.86"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "3.25"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.4147"
$ws.Range("H3").Value = 8

$a = $ws.Cells.Item(4, 1)
$a.Value = 2
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "163302"
$ws.Range("C4").Value = "大摩资源优选混合（LOF）"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.92"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "88.13"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "5.10"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0.2509"
$ws.Range("H4").Value = 3

$a = $ws.Cells.Item(5, 1)
$a.Value = 3
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "013627"
$ws.Range("C5").Value = "华夏周期驱动混合C"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "4.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "86.45"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "3.18"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.1549"
$ws.Range("H5").Value = 10

$a = $ws.Cells.Item(6, 1)
$a.Value = 4
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "011172"
$ws.Range("C6").Value = "广发利鑫灵活配置混合C"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.03"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "73.90"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2.14"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "0.1504"
$ws.Range("H6").Value = 10

$a = $ws.Cells.Item(7, 1)
$a.Value = 5
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "013626"
$ws.Range("C7").Value = "华夏周期驱动混合A"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.12"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "86.45"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "3.18"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.1310"
$ws.Range("H7").Value = 10

$a = $ws.Cells.Item(8, 1)
$a.Value = 6
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "011346"
$ws.Range("C8").Value = "淳厚鑫淳一年持有期混合"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.34"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "78.96"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2.39"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.0798"
$ws.Range("H8").Value = 9

$a = $ws.Cells.Item(9, 1)
$a.Value = 7
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "012454"
$ws.Range("C9").Value = "淳厚鑫悦混合A"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.82"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "85.29"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "2.85"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.0519"
$ws.Range("H9").Value = 5

$a = $ws.Cells.Item(10, 1)
$a.Value = 8
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "012455"
$ws.Range("C10").Value = "淳厚鑫悦混合C"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.57"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "85.29"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "2.85"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.0162"
$ws.Range("H10").Value = 5

$a = $ws.Cells.Item(11, 1)
$a.Value = 9
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "002584"
$ws.Range("C11").Value = "富安达长盈灵活配置混合A"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.10"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "85.18"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "3.15"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.0032"
$ws.Range("H11").Value = 8

$a = $ws.Cells.Item(12, 1)
$a.Value = 10
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "004795"
$ws.Range("C12").Value = "富荣福鑫灵活配置混合C"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.05"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "68.35"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "4.01"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.0020"
$ws.Range("H12").Value = 1

$a = $ws.Cells.Item(13, 1)
$a.Value = 11
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "004794"
$ws.Range("C13").Value = "富荣福鑫灵活配置混合A"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "68.35"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "4.01"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.0004"
$ws.Range("H13").Value = 1

$a = $ws.Cells.Item(14, 1)
$a.Value = 12
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "016214"
$ws.Range("C14").Value = "富安达长盈灵活配置混合C"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.01"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "85.18"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "3.15"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0.0003"
$ws.Range("H14").Value = 8

# ---------------------------------------------------------------
# 2) Update the "总计" summary sheet: add a 2022-Q4 row on top,
#    pushing the older quarters down by one row.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Rows.Item(2).ClearFormats()

$a2 = $total.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 13
$total.Range("D2").Value = 1.74
